$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new column A (existing columns A-D shift to B-E) ---
$ws.Columns("A:A").Insert()

# --- Insert a new row 2 (existing row 2 shifts to row 3) ---
$ws.Rows("2:2").Insert()

# --- Column A labels (TabName / CasesTab / FilesTab) ---
$ws.Range("A1").Value = "TabName"
$ws.Range("A2").Value = "CasesTab"
$ws.Range("A3").Value = "FilesTab"

# --- Row 3: FilesTab query (brand-new query text) ---
$filesQuery = @'
MATCH (f:file)
OPTIONAL MATCH (f)-[*]->(a:arm)-[:of_trial]->(ct:clinical_trial)
OPTIONAL MATCH (f)-[*]->(c:case)
OPTIONAL MATCH (f)-->(parent)
WITH f,a,ct,c,parent
WHERE a.arm_id IN ['Q']
WITH
    f, parent, c, a, ct,
    ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,
    toInteger(floor(log(f.file_size)/log(1024))) as i,
    2 as precision
WITH
    f, parent, c, a, ct,
    f.file_size /(1024^i) AS value,
    10^precision AS factor,
    units[i] as unit
WITH
    f, parent, c, a, ct, unit,
    round(factor * value)/factor AS size
RETURN DISTINCT
    f.file_name AS `File Name`,
    head(labels(parent)) as Association,
    f.file_description AS Description,
    f.file_format AS `File Format`,
    CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,
    ct.clinical_trial_designation AS `Trial Code`,
    a.arm_id AS Arm,
    c.case_id AS `Case ID`
'@
$ws.Range("B3").Value = $filesQuery

# --- Row 2: CasesTab query (rewritten query text) ---
$casesQuery = @'
MATCH (c:case)
 MATCH (c)-[:of_arm]->(a:arm)-[:of_trial]->(ct:clinical_trial)
 WHERE a.arm_id IN ['Q']
OPTIONAL MATCH (f:file)-[*]->(c)
RETURN DISTINCT
    c.case_id AS `Case ID`,
     ct.clinical_trial_designation AS `Trial Code`,
     a.arm_id AS Arm,
      a.arm_drug AS `Arm Treatment`,
c.disease AS Diagnosis,
  c.gender AS Gender,
    c.race AS Race,
    c.ethnicity AS Ethnicity
'@
$ws.Range("B2").Value = $casesQuery

# --- Shared StatsTab query used by both rows ---
$statsQuery = @'
MATCH (f:file)
OPTIONAL MATCH (f)-[*]->(a:arm)-[:of_trial]->(ct:clinical_trial)
OPTIONAL MATCH (f)-[*]->(c:case)
WITH f,a,ct,c
WHERE a.arm_id IN ['Q']
RETURN
    COUNT(DISTINCT ct.clinical_trial_designation) AS Trials,
    COUNT(DISTINCT c.case_id) AS Cases,
    COUNT(DISTINCT f) AS Files
'@
$ws.Range("C2").Value = $statsQuery
$ws.Range("C3").Value = $statsQuery

# --- File name columns (D/E) ---
$ws.Range("D2").Value = "TC01_Trials_Filter_TrialArm-Q_Neo4jData.xlsx"
$ws.Range("E2").Value = "TC01_Trials_Filter_TrialArm-Q_WebData.xlsx"

# --- Apply wrap-text style (existing style index 1) to query cells ---
$ws.Range("B2").WrapText = $true
$ws.Range("C2").WrapText = $true
$ws.Range("B3").WrapText = $true
$ws.Range("C3").WrapText = $true

# --- Row heights to match autofit wrapped text ---
$ws.Rows(2).RowHeight = 195
$ws.Rows(3).RowHeight = 409.5

# --- Column widths ---
$ws.Columns("A:A").ColumnWidth = 8.85546875
$ws.Columns("B:C").ColumnWidth = 75.85546875
$ws.Columns("D:D").ColumnWidth = 70.28515625
$ws.Columns("E:E").ColumnWidth = 28.5703125

# --- View settings ---
$ws.Range("B3").Select()
